# Simulator full-month coverage, persist logs, fix employees
# Updates the Jason Green 2026-01-12 weekly timesheet export:
#  - Row for 2026-01-12 becomes a PTO entry (6.5h @ 100 = 650)
#  - 2026-01-13 client renamed Evans -> Leixner/Smith (6.5h @ 100 = 650)
#  - 2026-01-14 client renamed Davis -> Hunter (7h @ 100 = 700)
#  - 2026-01-15 client renamed Funke -> O'Connor (6h @ 100 = 600)
#  - 2026-01-16 client renamed Field -> Varricchio (6h @ 100 = 600)
#  - Subtotal/grand total hours & dollars recomputed (32h / $3200)
#  - Employee ID reassigned to emp_qhpjptqm

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")
$ws2 = $wb.Worksheets.Item("Jason Schema")

# --- Weekly Timesheet sheet ---

# 2026-01-12 -> PTO day
$ws1.Range("B2").Value = "PTO"
$ws1.Range("C2").Value = 6.5
$ws1.Range("D2").Value = "PTO"
$ws1.Range("E2").Value = 100
$ws1.Range("F2").Value = 650

# 2026-01-13 -> Leixner/Smith
$ws1.Range("B3").Value = "Leixner/Smith"
$ws1.Range("C3").Value = 6.5
$ws1.Range("E3").Value = 100
$ws1.Range("F3").Value = 650

# 2026-01-14 -> Hunter
$ws1.Range("B4").Value = "Hunter"
$ws1.Range("C4").Value = 7
$ws1.Range("E4").Value = 100
$ws1.Range("F4").Value = 700

# 2026-01-15 -> O'Connor
$ws1.Range("B5").Value = "O'Connor"
$ws1.Range("C5").Value = 6
$ws1.Range("E5").Value = 100
$ws1.Range("F5").Value = 600

# 2026-01-16 -> Varricchio
$ws1.Range("B6").Value = "Varricchio"
$ws1.Range("C6").Value = 6
$ws1.Range("E6").Value = 100
$ws1.Range("F6").Value = 600

# Subtotal / totals row
$ws1.Range("C8").Value = 32
$ws1.Range("D8").Value = "Reg: 32 / OT: 0"
$ws1.Range("F8").Value = 3200
$ws1.Range("F11").Value = 3200
$ws1.Range("F13").Value = 3200

# --- Jason Schema sheet ---

$ws2.Range("D2").Value = "PTO"
$ws2.Range("E2").Value = 6.5
$ws2.Range("F2").Value = 100
$ws2.Range("G2").Value = 650
$ws2.Range("H2").Value = "PTO"
$ws2.Range("I2").Value = "PTO"

$ws2.Range("D3").Value = "Leixner/Smith"
$ws2.Range("E3").Value = 6.5
$ws2.Range("F3").Value = 100
$ws2.Range("G3").Value = 650

$ws2.Range("D4").Value = "Hunter"
$ws2.Range("E4").Value = 7
$ws2.Range("F4").Value = 100
$ws2.Range("G4").Value = 700

$ws2.Range("D5").Value = "O'Connor"
$ws2.Range("E5").Value = 6
$ws2.Range("F5").Value = 100
$ws2.Range("G5").Value = 600

$ws2.Range("D6").Value = "Varricchio"
$ws2.Range("E6").Value = 6
$ws2.Range("F6").Value = 100
$ws2.Range("G6").Value = 600

# Employee ID re-issued for all rows
$ws2.Range("B2:B6").Value = "emp_qhpjptqm"
